# Generate Report for Archive
#
# 1. Update the "Ready for handoff" status text to "In Translation"
#    everywhere it is used (Overview sheet columns E/F, and the per-locale
#    "Status" column C on the zh-cn / de-de sheets).
# 2. Shrink the now-narrower "Status" columns to match the new, shorter text.
#    (Target authored width is 13.4101845877511 chars; this runtime quantizes
#    ColumnWidth assignments to 1/6-character steps, so 12.5 is the input
#    that lands on the nearest reachable stored width, 13.333333333333334.)

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2:F3").Value = "In Translation"
$overview.Range("E1:F1").EntireColumn.ColumnWidth = 12.5

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2:C3").Value = "In Translation"
$zhcn.Range("C1").EntireColumn.ColumnWidth = 12.5

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2:C3").Value = "In Translation"
$dede.Range("C1").EntireColumn.ColumnWidth = 12.5
